$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C19").Value = 'Plantations'

$ws.Range("E21").Value = 'Carracks'

$ws.Range("E23").Value = 'Coastal Batteries'

$ws.Range("C25").Value = 'Textile Manufactory'

$ws.Range("E26").Value = 'Colonial Charters'

$ws.Range("C28").Value = 'Modern Metallurgy'

$ws.Range("C29").Value = 'Scientific Enquiry'
$ws.Range("E29").Value = 'Screw Propeller'

$ws.Range("E30").Value = 'Merchantman'

$ws.Range("C33").Value = 'University'
$ws.Range("E33").Value = 'Threedeckers'

$ws.Range("E34").Value = 'Naval Proffessionalization'

$ws.Range("E35").Value = 'Steam Turbine'

$ws.Range("E36").Value = 'Merchant Vessels'

$ws.Range("E37").Value = 'All-or-Nothing Armor Scheme'

$ws.Range("E38").Value = 'Joint Stock Companies'

$ws.Range("E39").Value = 'Public Punishments'
$ws.Range("G39").Value = 'Breech-loading'

$ws.Range("E42").Value = 'Cargo Ships'

$ws.Range("C44").Value = 'Railroads'

$ws.Range("G45").Value = 'Lile Rifle'
$ws.Range("H45").Value = '+0.5 inf fire, +50% supply limit, +25% flanking range, new inf'

$ws.Range("C46").Value = 'Mechanized Mining'
$ws.Range("E46").Value = 'Market Regulation'
$ws.Range("G46").Value = 'Rifled Artillery'

$ws.Range("C47").Value = 'Tractors'

$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("C47").Select()
